$d = $word.ActiveDocument

# The document currently has a single empty paragraph. Put the first
# line of text into it, then append the remaining lines as new
# paragraphs before the end of the document (end mark).

$para1 = $d.Paragraphs.Item(1)
$para1.Range.Text = "Présenter l’intérêt des visualisations clés de l’application"

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)
$endRange.InsertAfter("Présenter les fonctionnalités majeures de l’application")

$endRange2 = $d.Content
$endRange2.Collapse(0)
$endRange2.InsertParagraphAfter()
$endRange2.Collapse(0)
$endRange2.InsertAfter("/ !\ Attention à la mise en forme et à la rédaction du docu / !\")
